$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three "description" cells (column G) with new wording.
$ws.Cells.Item(6, 7).Value = "Données démographiques et géographiques par région"
$ws.Cells.Item(5, 7).Value = "Évolution des températures mensuelles"
$ws.Cells.Item(17, 7).Value = "Analyse des prix des biens immobiliers dans les grandes villes"

# Move the active selection on the sheet to D24 (matches the saved view state).
$ws.Range("D24").Select()
